$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slovenia Prva Liga")

# The underlying data feed re-sent this league's fixtures with two pairs of
# rows transposed (the match that used to sit in one row now sits in the
# other, and vice-versa). Column A (the running id) stays put; everything
# from column B (match id) through column AD (closing odds) swaps between
# the two rows for each pair.

function Swap-Rows($ws, $row1, $row2, $firstCol, $lastCol) {
    $range1 = $ws.Range($ws.Cells.Item($row1, $firstCol), $ws.Cells.Item($row1, $lastCol))
    $range2 = $ws.Range($ws.Cells.Item($row2, $firstCol), $ws.Cells.Item($row2, $lastCol))

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

# Columns B (2) through AD (30)
$firstCol = 2
$lastCol = 30

# Rows 9 and 10 (matches 6814328 / 6814330) swap places
Swap-Rows $ws 9 10 $firstCol $lastCol

# Rows 175 and 176 (matches 7133777 / 7124153) swap places
Swap-Rows $ws 175 176 $firstCol $lastCol
